$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("D2").Value = "Comments for even level"
$ws.Range("D3").Value = "handover count per case"
$ws.Range("D4").Value = "bit stating whether it is rework or not"
$ws.Range("D5").Value = "not included in the event_consolidated"

$ws.Range("D19").Select()
